$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the three runs (split around a gramStart/gramEnd
# proofErr pair) that make up the "counter-example" sentence into a
# single run, dropping the proofErr markers. The visible text is
# unchanged; only the run/proofErr structure collapses.
# ---------------------------------------------------------------------
$target = "The counter-example I have chosen is the Golden Rule, " + [char]0x201C + "Do unto others as you would have them do unto you." + [char]0x201D

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pText = $p.Range.Text.TrimEnd([char]13)
    if ($pText -eq $target) {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $txt = $r.Text
        $r.Delete()
        $r2 = $d.Range($p.Range.Start, $p.Range.Start)
        $r2.InsertAfter($txt)
        $found = $true
        break
    }
}
Write-Output "counter-example paragraph collapsed: $found"

# ---------------------------------------------------------------------
# Change 2: after the final paragraph of the document, append three new
# paragraphs:
#   - an empty ListParagraph (ind left=1440)
#   - a ListParagraph numbered item (ilvl 0, numId 3) containing " "
#   - a ListParagraph numbered item (ilvl 1, numId 3) containing the new
#     "As a search problem..." text
#
# Insert them in one shot as raw WordprocessingML via InsertXML at the
# end of the current last (non-empty) paragraph, so Word splits a new
# paragraph off cleanly instead of merging into an existing empty one.
# ---------------------------------------------------------------------
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newPara1 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p>'

$newPara2 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$newPara3Text = "As a search problem, I would define the states as possible steps in a solution to the problem. The start state would be a tuple with the first index and the " + [char]0x201C + "-BEGIN-" + [char]0x201D + " string. Each subsequent state would be the previous word and the index of the next letter. The end state would be whenever the state[0] value is equal to the length of the query. Costs would be determined with a bigram cost function with parameters (previous word, new proposed word). "
$newPara3 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $newPara3Text + '</w:t></w:r></w:p>'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
[void]$insertionPoint.InsertXML($newPara1 + $newPara2 + $newPara3)

Write-Output "final paragraph count: $($d.Paragraphs.Count)"
